$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Catálogo (catalog.html)", "Muestra todos los productos agrarios disponibles en la plataforma. Este se muestra al pulsar sobre catalogo (mostrará todo), al realizar una búsqueda por nombre de producto o por ubicación (mostrará solo las coincidencias) o al pulsar sobre una categoría (mostrará solo dicha categoría), además permitirá acceder a la información`n de cada producto o añadirlo a la cesta directamente."),
    @("Ofertas (offers.html)", "Muestra el listado de ofertas disponibles en ese momento."),
    @("Categorías (category.html)", "Muestra todas las categorías de productos que hay en la plataforma."),
    @("Producto (product.html)", "Muestra la información de un producto y permitirá añadirlo a la cesta.")
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A15:B15").Copy()
    $ws.Range("A" + $r + ":B" + $r).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Rows.Item($r).EntireRow.AutoFit()
}
